$d = $word.ActiveDocument

function Rename-InlineShapeImage($inlineShape, $newName) {
    $floatShape = $inlineShape.ConvertToShape()
    $floatShape.Name = $newName
    return $floatShape.ConvertToInlineShape()
}

$sec = $d.Sections(1)

# Footer 1: Pearson logo (image2.png -> image1.png)
$ftr1 = $sec.Footers(1)
$shp1 = $ftr1.Range.InlineShapes(1)
Rename-InlineShapeImage $shp1 "image1.png" | Out-Null

# Footer 2: Pearson logo (image2.png -> image1.png)
$ftr2 = $sec.Footers(2)
$shp2 = $ftr2.Range.InlineShapes(1)
Rename-InlineShapeImage $shp2 "image1.png" | Out-Null

# Header 2: BTEC logo (image1.jpg -> image2.jpg)
$hdr2 = $sec.Headers(2)
$shp3 = $hdr2.Range.InlineShapes(1)
Rename-InlineShapeImage $shp3 "image2.jpg" | Out-Null

Write-Output "Renamed 3 inline shape images"
